$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap date and volumen values between row 2 and row 5
$ws.Range("D2").Value = 44691
$ws.Range("J2").Value = 100

$ws.Range("D5").Value = 44692
$ws.Range("J5").Value = 120
